$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data appended to the bottom of the sheet (rows 199-210),
# continuing the existing series in columns A (index) and B (value).
$newData = @(
    @(197, [double]"2.775557561562891E-17"),
    @(198, [double]"-2.775557561562891E-17"),
    @(199, [double]"-6.245004513516506E-17"),
    @(200, [double]"6.630498619289129E-17"),
    @(201, [double]"1.908195823574488E-17"),
    @(202, [double]"-2.478176394252582E-17"),
    @(203, [double]"2.775557561562891E-17"),
    @(204, [double]"2.081668171172168E-18"),
    @(205, [double]"8.673617379884035E-19"),
    @(206, [double]"5.406554833461049E-17"),
    @(207, [double]"0"),
    @(208, [double]"0")
)

$lastRow = 198
$styleSource = $ws.Range("A$lastRow")

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $lastRow + 1 + $i
    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)

    $aCell.Value2 = $newData[$i][0]
    $bCell.Value2 = $newData[$i][1]

    # Match the formatting used by the existing index column (bold, bordered,
    # centered) by copying the style from the last pre-existing row.
    $styleSource.Copy()
    $aCell.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
